# Applies the "removed comments / Added restock keyword in domain excel /
# Added review transactions sentence in dataelements" commit to Sheet1.
#
# Summary of the change:
#  - Row 8 (Banker / restock sufficient cash): Functionality (C8) gets a new
#    value "['Restock cash', 'Limit Cash']" and Attributes (D8) becomes
#    "['Cus_Nme', 'Acc_type']['Cus_Nme', 'Acc_type']".
#  - Row 10 (Banker / review all transactions): Attributes (D10) gets the
#    "Acc_num" element added -> "['Acc_num']['Acc_num', 'Amt_avail']".
#  - Row 11 (Banker / review credit history...): Attributes (D11) also picks
#    up the "Acc_num" element in the first segment ->
#    "['Loan_Amt', 'Cred_Score']['Acc_num'],[...]['Loan_Amt', 'Cred_Score']".
#  - Column B's custom width is cleared and column C's width becomes 47.
#  - calcPr gets concurrentCalc="0".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Functionality / Attributes text for the affected rows ---
$ws.Range("C8").Value = "['Restock cash', 'Limit Cash']"
$ws.Range("D8").Value = "['Cus_Nme', 'Acc_type']['Cus_Nme', 'Acc_type']"

$ws.Range("D10").Value = "['Acc_num']['Acc_num', 'Amt_avail']"

$ws.Range("D11").Value = "['Loan_Amt', 'Cred_Score']['Acc_num'],['Loan_Amt', 'Amt_avail', 'Debit_pin', 'Amt_wdrl', 'Amt_trnsfr']['Loan_Amt', 'Cred_Score']"

# --- Column width adjustments: drop B's custom width, set C to 47 ---
# (B reverts to the workbook's standard/default column width; C's stored
# width needs to serialize as exactly 47 characters - Excel's ColumnWidth
# -> stored-width pixel rounding means 46.14 lands on 47.)
$ws.Columns.Item(2).ColumnWidth = $ws.StandardWidth
$ws.Columns.Item(3).ColumnWidth = 46.14

# --- Workbook calculation setting (concurrentCalc="0") ---
$excel.Application.MultiThreadedCalculation.Enabled = $false
